$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 72125
$ws.Range("J3").Value = 72125
$ws.Range("L3").Value = 72125
$ws.Range("N3").Value = -72353
$ws.Range("H40").Value = 6428.385
$ws.Range("I40").Value = 3755
$ws.Range("J40").Value = 8719.857
$ws.Range("K40").Value = 3755
$ws.Range("L40").Value = 8719.857
$ws.Range("M40").Value = -3580
$ws.Range("N40").Value = -9069.857
$ws.Range("H62").Value = 4967.3335
$ws.Range("I62").Value = 3951
$ws.Range("K62").Value = 3951
$ws.Range("M62").Value = -3327
$ws.Range("H65").Value = 4967.3335
$ws.Range("I65").Value = 3951
$ws.Range("K65").Value = 19755
$ws.Range("M65").Value = -16635
$ws.Range("H102").Value = 72125
$ws.Range("J102").Value = 72125
$ws.Range("L102").Value = 72125
$ws.Range("N102").Value = -78615
$ws.Range("H132").Value = 13306.24
$ws.Range("I132").Value = 2174.476
$ws.Range("K132").Value = 6523.428
$ws.Range("M132").Value = -3993.428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16400255
$ws.Range("I32").Value = 17550816
$ws.Range("J32").Value = 4747.25
$ws.Range("K32").Value = 17550816
$ws.Range("L32").Value = 4747.25
$ws.Range("M32").Value = -17550529
$ws.Range("N32").Value = -5321.25
$ws.Range("H61").Value = 3924.875
$ws.Range("I61").Value = 4199.857
$ws.Range("K61").Value = 4199.857
$ws.Range("M61").Value = -3987.857
$ws.Range("H74").Value = 2851.5
$ws.Range("I74").Value = 2886.3076
$ws.Range("K74").Value = 2886.3076
$ws.Range("M74").Value = -2012.3076
$ws.Range("H77").Value = 2851.5
$ws.Range("I77").Value = 2886.3076
$ws.Range("K77").Value = 14431.538
$ws.Range("M77").Value = -10063.538
$ws.Range("H125").Value = 70000
$ws.Range("J125").Value = 70000
$ws.Range("L125").Value = 70000
$ws.Range("N125").Value = -79840
$ws.Range("H132").Value = 2265.7693
$ws.Range("I132").Value = 2199.9583
$ws.Range("J132").Value = 3055.5
$ws.Range("K132").Value = 6599.874899999999
$ws.Range("L132").Value = 9166.5
$ws.Range("M132").Value = -4069.874899999999
$ws.Range("N132").Value = -14226.5
$ws.Range("H136").Value = 3924.875
$ws.Range("I136").Value = 4199.857
$ws.Range("K136").Value = 12599.571
$ws.Range("M136").Value = -10049.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3543.5293
$ws.Range("I20").Value = 1232.3
$ws.Range("J20").Value = 6845.2856
$ws.Range("K20").Value = 1232.3
$ws.Range("L20").Value = 6845.2856
$ws.Range("M20").Value = -985.3
$ws.Range("N20").Value = -7339.2856
$ws.Range("H86").Value = 4515.6
$ws.Range("I86").Value = 3338.1667
$ws.Range("K86").Value = 3338.1667
$ws.Range("M86").Value = -2215.1667
$ws.Range("H89").Value = 4515.6
$ws.Range("I89").Value = 3338.1667
$ws.Range("K89").Value = 16690.8335
$ws.Range("M89").Value = -11074.8335
$ws.Range("H94").Value = 3567.5557
$ws.Range("I94").Value = 4231.615
$ws.Range("J94").Value = 1841
$ws.Range("K94").Value = 4231.615
$ws.Range("L94").Value = 1841
$ws.Range("M94").Value = -3780.615
$ws.Range("N94").Value = -2743
$ws.Range("H99").Value = 67703.164
$ws.Range("I99").Value = 100577.25
$ws.Range("K99").Value = 100577.25
$ws.Range("M99").Value = -99079.25
$ws.Range("H141").Value = 48000
$ws.Range("J141").Value = 48000
$ws.Range("L141").Value = 48000
$ws.Range("N141").Value = -58360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1905.0731
$ws.Range("I31").Value = 1799.3429
$ws.Range("K31").Value = 1799.3429
$ws.Range("M31").Value = -1504.3429
$ws.Range("H34").Value = 1905.0731
$ws.Range("I34").Value = 1799.3429
$ws.Range("K34").Value = 1799.3429
$ws.Range("M34").Value = -1597.3429
$ws.Range("H58").Value = 1973.2424
$ws.Range("J58").Value = 4186.9
$ws.Range("L58").Value = 4186.9
$ws.Range("N58").Value = -4592.9
$ws.Range("H109").Value = 17160.188
$ws.Range("J109").Value = 17160.188
$ws.Range("L109").Value = 17160.188
$ws.Range("N109").Value = -19240.188
$ws.Range("H132").Value = 1813.5217
$ws.Range("I132").Value = 1759.591
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 5278.772999999999
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -2748.772999999999
$ws.Range("N132").Value = -14060
$ws.Range("H134").Value = 1932.8649
$ws.Range("I134").Value = 1961.8823
$ws.Range("J134").Value = 1604
$ws.Range("K134").Value = 5885.6469
$ws.Range("L134").Value = 4812
$ws.Range("M134").Value = -3350.6469
$ws.Range("N134").Value = -9882
$ws.Range("H136").Value = 1973.2424
$ws.Range("J136").Value = 4186.9
$ws.Range("L136").Value = 12560.7
$ws.Range("N136").Value = -17660.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 840
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 840
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2520
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -5016
$ws.Range("H107").Value = 756
$ws.Range("J107").Value = 844
$ws.Range("L107").Value = 2532
$ws.Range("N107").Value = -6372
$ws.Range("H113").Value = 1493.6923
$ws.Range("I113").Value = 1331.4286
$ws.Range("J113").Value = 1683
$ws.Range("K113").Value = 3994.2858
$ws.Range("L113").Value = 5049
$ws.Range("M113").Value = -1824.2858
$ws.Range("N113").Value = -9389
$ws.Range("H114").Value = 4379.8
$ws.Range("H129").Value = 1730.2667
$ws.Range("J129").Value = 2311.875
$ws.Range("L129").Value = 6935.625
$ws.Range("N129").Value = -16935.625
$ws.Range("H131").Value = 3376.5757
$ws.Range("J131").Value = 3700.4167
$ws.Range("L131").Value = 11101.2501
$ws.Range("N131").Value = -21181.2501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3917.6667
$ws.Range("I43").Value = 703.2
$ws.Range("K43").Value = 703.2
$ws.Range("M43").Value = -552.2
$ws.Range("H80").Value = 98570.62
$ws.Range("I80").Value = 149570.5
$ws.Range("K80").Value = 149570.5
$ws.Range("M80").Value = -148572.5
$ws.Range("H83").Value = 98570.62
$ws.Range("I83").Value = 149570.5
$ws.Range("K83").Value = 747852.5
$ws.Range("M83").Value = -742860.5
$ws.Range("H122").Value = 7584.357
$ws.Range("I122").Value = 9294.5
$ws.Range("J122").Value = 6301.75
$ws.Range("K122").Value = 27883.5
$ws.Range("L122").Value = 18905.25
$ws.Range("M122").Value = -25433.5
$ws.Range("N122").Value = -23805.25
$ws.Range("H132").Value = 2321.3333
$ws.Range("I132").Value = 1772.8889
$ws.Range("K132").Value = 5318.6667
$ws.Range("M132").Value = -2788.6667
$ws.Range("H135").Value = 59991.85
$ws.Range("J135").Value = 59991.85
$ws.Range("L135").Value = 59991.85
$ws.Range("N135").Value = -70131.85000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2013.8334
$ws.Range("I16").Value = 2316.6
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 2316.6
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -2146.6
$ws.Range("N16").Value = -840
$ws.Range("H40").Value = 10044.929
$ws.Range("I40").Value = 10872.833
$ws.Range("K40").Value = 10872.833
$ws.Range("M40").Value = -10736.833
$ws.Range("H61").Value = 1838.5927
$ws.Range("I61").Value = 1393.8636
$ws.Range("K61").Value = 1393.8636
$ws.Range("M61").Value = -1191.8636
$ws.Range("H113").Value = 1838.5927
$ws.Range("I113").Value = 1393.8636
$ws.Range("K113").Value = 1393.8636
$ws.Range("M113").Value = 776.1364000000001
$ws.Range("H132").Value = 6282.125
$ws.Range("I132").Value = 3620.4
$ws.Range("J132").Value = 10718.333
$ws.Range("K132").Value = 10861.2
$ws.Range("L132").Value = 32154.999
$ws.Range("M132").Value = -8331.200000000001
$ws.Range("N132").Value = -37214.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()
$ws.Range("H96").Value = 39840.785
$ws.Range("J96").Value = 4828
$ws.Range("L96").Value = 4828
$ws.Range("N96").Value = -7574
$ws.Range("H132").Value = 2919.6
$ws.Range("I132").Value = 2946.4285
$ws.Range("K132").Value = 8839.2855
$ws.Range("M132").Value = -6309.2855
